$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = -0
$ws.Cells.Item(2, 2).Value = -0.09990638128057558
$ws.Cells.Item(2, 3).Value = -0
$ws.Cells.Item(2, 4).Value = 0.2120382236287753
$ws.Cells.Item(2, 5).Value = 0.01537492257737358
$ws.Cells.Item(2, 6).Value = -0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 10).Value = -0
$ws.Cells.Item(2, 11).Value = -0.07541003767991721
$ws.Cells.Item(2, 12).Value = -0
$ws.Cells.Item(2, 13).Value = 0.1924694857085359
$ws.Cells.Item(2, 14).Value = 0.02980939654462182
$ws.Cells.Item(2, 18).Value = -0
$ws.Cells.Item(2, 19).Value = -0
$ws.Cells.Item(2, 20).Value = -0.112993313119286
$ws.Cells.Item(2, 21).Value = -0
$ws.Cells.Item(2, 22).Value = 0.01741763980473633
$ws.Cells.Item(2, 23).Value = -0.05480521067106651
$ws.Cells.Item(2, 25).Value = 0
$ws.Cells.Item(2, 26).Value = 0
$ws.Cells.Item(2, 28).Value = 0
$ws.Cells.Item(2, 29).Value = -0.04763560623632488
$ws.Cells.Item(2, 30).Value = 0
$ws.Cells.Item(2, 31).Value = -0.001231494748786342
$ws.Cells.Item(2, 32).Value = 0.02180049073613898
$ws.Cells.Item(2, 33).Value = -0
$ws.Cells.Item(2, 35).Value = -0
$ws.Cells.Item(2, 36).Value = 0
$ws.Cells.Item(2, 37).Value = -0
$ws.Cells.Item(2, 38).Value = -0.017011889677344
$ws.Cells.Item(2, 40).Value = 0.01655520015105
$ws.Cells.Item(2, 41).Value = 0.08616768995696322
$ws.Cells.Item(2, 43).Value = 0
$ws.Cells.Item(2, 45).Value = 0
$ws.Cells.Item(2, 46).Value = 0
$ws.Cells.Item(2, 47).Value = 0.008540910523613702
$ws.Cells.Item(2, 49).Value = 0.02384383279224323
$ws.Cells.Item(2, 50).Value = 0.002653279452716468
$ws.Cells.Item(2, 51).Value = -0
$ws.Cells.Item(2, 55).Value = -0
$ws.Cells.Item(2, 56).Value = -0.05609398520719405
$ws.Cells.Item(2, 58).Value = 0.0901994727725183
$ws.Cells.Item(2, 59).Value = 0.03183116576526138
$ws.Cells.Item(2, 61).Value = 0
$ws.Cells.Item(2, 62).Value = -0
$ws.Cells.Item(2, 63).Value = -0
$ws.Cells.Item(2, 64).Value = 0
$ws.Cells.Item(2, 65).Value = -0.01740503018763945
$ws.Cells.Item(2, 67).Value = -0.04817202116933253
$ws.Cells.Item(2, 68).Value = -0.0458225377765082
$ws.Cells.Item(2, 73).Value = 0
$ws.Cells.Item(2, 74).Value = -0.03495484539737654
$ws.Cells.Item(2, 76).Value = 0.0121282688106168
$ws.Cells.Item(2, 77).Value = -0.01135934596789276
$ws.Cells.Item(2, 78).Value = -0
$ws.Cells.Item(2, 80).Value = 0
$ws.Cells.Item(2, 81).Value = 0
$ws.Cells.Item(2, 82).Value = 0
$ws.Cells.Item(2, 83).Value = 0.02622384571045857
$ws.Cells.Item(2, 85).Value = -0.01464824478628676
$ws.Cells.Item(2, 86).Value = -0.003218958377605594
$ws.Cells.Item(2, 88).Value = -0
$ws.Cells.Item(2, 91).Value = -0
$ws.Cells.Item(2, 92).Value = -0.02910295133324822
$ws.Cells.Item(2, 93).Value = -0
$ws.Cells.Item(2, 94).Value = 0.06160945795060888
$ws.Cells.Item(2, 95).Value = 0.06346280035528887
$ws.Cells.Item(2, 98).Value = 0
$ws.Cells.Item(2, 99).Value = -0
$ws.Cells.Item(2, 100).Value = -0
$ws.Cells.Item(2, 101).Value = 0.03159391713659575
$ws.Cells.Item(2, 103).Value = -0.02389457700439125
$ws.Cells.Item(2, 104).Value = -0.002725825012946287
$ws.Cells.Item(2, 109).Value = -0
$ws.Cells.Item(2, 110).Value = 0.03764982877090165
$ws.Cells.Item(2, 112).Value = 0.0241140576307954
$ws.Cells.Item(2, 113).Value = 0.02752073542778462
$ws.Cells.Item(2, 114).Value = 0
$ws.Cells.Item(2, 115).Value = -0
$ws.Cells.Item(2, 116).Value = -0
$ws.Cells.Item(2, 118).Value = 0
$ws.Cells.Item(2, 119).Value = 0.01583361300337519
$ws.Cells.Item(2, 120).Value = -0
$ws.Cells.Item(2, 121).Value = 0.0730609548709788
$ws.Cells.Item(2, 122).Value = 0.04441029524461618
$ws.Cells.Item(2, 124).Value = 0
$ws.Cells.Item(2, 125).Value = -0
$ws.Cells.Item(2, 126).Value = -0
$ws.Cells.Item(2, 127).Value = 0
$ws.Cells.Item(2, 128).Value = -0.02475748232067152
$ws.Cells.Item(2, 129).Value = -0
$ws.Cells.Item(2, 130).Value = 0.009866557565529403
$ws.Cells.Item(2, 131).Value = -0.0194466426991889
$ws.Cells.Item(2, 132).Value = 0
$ws.Cells.Item(2, 136).Value = -0
$ws.Cells.Item(2, 137).Value = 0.02151702844159892
$ws.Cells.Item(2, 139).Value = 0.08375347685187916
$ws.Cells.Item(2, 140).Value = -0.06249611160263986
$ws.Cells.Item(2, 142).Value = -0
$ws.Cells.Item(2, 144).Value = -0
$ws.Cells.Item(2, 145).Value = 0
$ws.Cells.Item(2, 146).Value = 0.02732689672173665
$ws.Cells.Item(2, 147).Value = 0
$ws.Cells.Item(2, 148).Value = -0.05271352491748541
$ws.Cells.Item(2, 149).Value = 0.0208501322478133
$ws.Cells.Item(2, 154).Value = 0
$ws.Cells.Item(2, 155).Value = 0.02384157451425638
$ws.Cells.Item(2, 156).Value = 0
$ws.Cells.Item(2, 157).Value = -0.01968171241820138
$ws.Cells.Item(2, 158).Value = -0.001258152351482244
$ws.Cells.Item(2, 159).Value = -0
$ws.Cells.Item(2, 160).Value = -0
$ws.Cells.Item(2, 162).Value = -0
$ws.Cells.Item(2, 163).Value = -0
$ws.Cells.Item(2, 164).Value = -0.00624867786378639
$ws.Cells.Item(2, 166).Value = -0.001504019357684276
$ws.Cells.Item(2, 167).Value = -0.04666068744010858
$ws.Cells.Item(2, 172).Value = -0
$ws.Cells.Item(2, 173).Value = -0.02887421200155625
$ws.Cells.Item(2, 174).Value = -0
$ws.Cells.Item(2, 175).Value = 0.006761150360827417
$ws.Cells.Item(2, 176).Value = 0.001189711286788366
$ws.Cells.Item(2, 177).Value = 0
$ws.Cells.Item(2, 179).Value = -0
$ws.Cells.Item(2, 180).Value = 0
$ws.Cells.Item(2, 181).Value = 0
$ws.Cells.Item(2, 182).Value = -0.05896639299117197
$ws.Cells.Item(2, 184).Value = 0.02240035335631822
$ws.Cells.Item(2, 186).Value = 0
$ws.Cells.Item(2, 187).Value = -0
$ws.Cells.Item(2, 188).Value = 0
